# Calibrate pixel->cm table: neutral points for ToF vs camera now differ,
# so the calibration series is replaced/extended (rows 2-53) and the
# camera curve continues out to A=180 / B=145.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New calibration data (A: cm, B: pixels) for rows 2..53.
$data = @(
  @(2,52),
  @(3,56),
  @(4,60),
  @(5,64),
  @(6,67),
  @(7,70),
  @(8,72),
  @(9,75),
  @(10,77),
  @(11,80),
  @(12,82),
  @(13,84),
  @(14,87),
  @(15,89),
  @(16,91),
  @(17,92),
  @(18,94),
  @(19,96),
  @(20,97),
  @(21,98),
  @(25,103.5),
  @(30,109),
  @(35,114),
  @(40,118),
  @(45,121),
  @(50,124),
  @(55,127),
  @(60,129),
  @(65,131),
  @(70,133),
  @(75,133.5),
  @(80,134.5),
  @(85,136),
  @(90,137),
  @(95,137.5),
  @(100,138.5),
  @(105,139.5),
  @(110,140),
  @(115,141),
  @(120,141),
  @(125,141),
  @(130,141.5),
  @(135,142),
  @(140,142.5),
  @(145,143),
  @(150,143.5),
  @(155,143.5),
  @(160,144.5),
  @(165,145),
  @(170,145),
  @(175,145),
  @(180,145)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $r = $i + 2
  $ws.Cells.Item($r, 1).Value = $data[$i][0]
  $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

# Selection now spans the whole filled table (A1:B53).
$ws.Range("A1:B53").Select()
